$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 666.6
$ws.Range("I33").Value = 552.2727
$ws.Range("K33").Value = 552.2727
$ws.Range("M33").Value = -323.2727

$ws.Range("H55").Value = 266
$ws.Range("I55").Value = 251.42857
$ws.Range("J55").Value = 300
$ws.Range("K55").Value = 251.42857
$ws.Range("L55").Value = 300
$ws.Range("M55").Value = -37.42857000000001
$ws.Range("N55").Value = -728

$ws.Range("H64").Value = 4279.4165
$ws.Range("I64").Value = 4039.8
$ws.Range("J64").Value = 5477.5
$ws.Range("K64").Value = 4039.8
$ws.Range("L64").Value = 5477.5
$ws.Range("M64").Value = -3791.8
$ws.Range("N64").Value = -5973.5

$ws.Range("H67").Value = 4279.4165
$ws.Range("I67").Value = 4039.8
$ws.Range("J67").Value = 5477.5
$ws.Range("K67").Value = 4039.8
$ws.Range("L67").Value = 5477.5
$ws.Range("M67").Value = -3181.8
$ws.Range("N67").Value = -7193.5

$ws.Range("H76").Value = 8086.7036
$ws.Range("J76").Value = 6120
$ws.Range("L76").Value = 6120
$ws.Range("N76").Value = -6750

$ws.Range("H79").Value = 8086.7036
$ws.Range("J79").Value = 6120
$ws.Range("L79").Value = 6120
$ws.Range("N79").Value = -8304

$ws.Range("H113").Value = 3425.9644
$ws.Range("I113").Value = 2455.1333
$ws.Range("J113").Value = 4546.154
$ws.Range("K113").Value = 2455.1333
$ws.Range("L113").Value = 4546.154
$ws.Range("M113").Value = 798.8667
$ws.Range("N113").Value = -11054.154

$ws.Range("H132").Value = 3023.2837
$ws.Range("I132").Value = 1750.2656
$ws.Range("J132").Value = 11170.6
$ws.Range("K132").Value = 5250.7968
$ws.Range("L132").Value = 33511.8
$ws.Range("M132").Value = -2720.7968
$ws.Range("N132").Value = -38571.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2011.7059
$ws.Range("I2").Value = 2727.9
$ws.Range("J2").Value = 988.5714
$ws.Range("K2").Value = 2727.9
$ws.Range("L2").Value = 988.5714
$ws.Range("M2").Value = -2614.9
$ws.Range("N2").Value = -1214.5714

$ws.Range("H32").Value = 27677.264
$ws.Range("I32").Value = 14047.615
$ws.Range("K32").Value = 14047.615
$ws.Range("M32").Value = -13760.615

$ws.Range("H34").Value = 15000
$ws.Range("J34").Value = 15000
$ws.Range("L34").Value = 15000
$ws.Range("N34").Value = -15542

$ws.Range("H63").Value = 14360
$ws.Range("I63").Value = 14360
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 14360
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -13674

$ws.Range("H66").Value = 14360
$ws.Range("I66").Value = 14360
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 71800
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -68368

$ws.Range("H88").Value = 1953.3334
$ws.Range("I88").Value = 1973.5294
$ws.Range("J88").Value = 1890.909
$ws.Range("K88").Value = 1973.5294
$ws.Range("L88").Value = 1890.909
$ws.Range("M88").Value = -1567.5294
$ws.Range("N88").Value = -2702.909

$ws.Range("H91").Value = 1953.3334
$ws.Range("I91").Value = 1973.5294
$ws.Range("J91").Value = 1890.909
$ws.Range("K91").Value = 1973.5294
$ws.Range("L91").Value = 1890.909
$ws.Range("M91").Value = -569.5293999999999
$ws.Range("N91").Value = -4698.909

$ws.Range("H97").Value = 2400
$ws.Range("I97").Value = 2542
$ws.Range("J97").Value = 1926.6666
$ws.Range("K97").Value = 2542
$ws.Range("L97").Value = 1926.6666
$ws.Range("M97").Value = -2046
$ws.Range("N97").Value = -2918.6666

$ws.Range("H116").Value = 2011.7059
$ws.Range("I116").Value = 2727.9
$ws.Range("J116").Value = 988.5714
$ws.Range("K116").Value = 2727.9
$ws.Range("L116").Value = 988.5714
$ws.Range("M116").Value = -433.9000000000001
$ws.Range("N116").Value = -5576.5714

$ws.Range("H122").Value = 1783.6511
$ws.Range("I122").Value = 1878.8611
$ws.Range("J122").Value = 1294
$ws.Range("K122").Value = 5636.5833
$ws.Range("L122").Value = 3882
$ws.Range("M122").Value = -3186.5833
$ws.Range("N122").Value = -8782

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2011.7059
$ws.Range("I3").Value = 2727.9
$ws.Range("J3").Value = 988.5714
$ws.Range("K3").Value = 2727.9
$ws.Range("L3").Value = 988.5714
$ws.Range("M3").Value = -2613.9
$ws.Range("N3").Value = -1216.5714

$ws.Range("H80").Value = 749.55554
$ws.Range("J80").Value = 749.55554
$ws.Range("L80").Value = 749.55554
$ws.Range("N80").Value = -2745.55554

$ws.Range("H83").Value = 749.55554
$ws.Range("J83").Value = 749.55554
$ws.Range("L83").Value = 3747.7777
$ws.Range("N83").Value = -13731.7777

$ws.Range("H86").Value = 4638.9707
$ws.Range("I86").Value = 3949
$ws.Range("J86").Value = 6081.636
$ws.Range("K86").Value = 3949
$ws.Range("L86").Value = 6081.636
$ws.Range("M86").Value = -2826
$ws.Range("N86").Value = -8327.636

$ws.Range("H89").Value = 4638.9707
$ws.Range("I89").Value = 3949
$ws.Range("J89").Value = 6081.636
$ws.Range("K89").Value = 19745
$ws.Range("L89").Value = 30408.18
$ws.Range("M89").Value = -14129
$ws.Range("N89").Value = -41640.18

$ws.Range("H105").Value = 2299.049
$ws.Range("I105").Value = 2113.125
$ws.Range("J105").Value = 2985.5386
$ws.Range("K105").Value = 2113.125
$ws.Range("L105").Value = 2985.5386
$ws.Range("M105").Value = -366.125
$ws.Range("N105").Value = -6479.5386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4175
$ws.Range("I62").Value = 2300
$ws.Range("J62").Value = 4800
$ws.Range("K62").Value = 2300
$ws.Range("L62").Value = 4800
$ws.Range("M62").Value = -1676
$ws.Range("N62").Value = -6048

$ws.Range("H65").Value = 4175
$ws.Range("I65").Value = 2300
$ws.Range("J65").Value = 4800
$ws.Range("K65").Value = 11500
$ws.Range("L65").Value = 24000
$ws.Range("M65").Value = -8380
$ws.Range("N65").Value = -30240

$ws.Range("H132").Value = 17244156
$ws.Range("I132").Value = 38462696
$ws.Range("J132").Value = 4093.625
$ws.Range("K132").Value = 115388088
$ws.Range("L132").Value = 12280.875
$ws.Range("M132").Value = -115385558
$ws.Range("N132").Value = -17340.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 928.5714
$ws.Range("I61").Value = 533.3333
$ws.Range("J61").Value = 1225
$ws.Range("K61").Value = 1599.9999
$ws.Range("L61").Value = 3675
$ws.Range("M61").Value = -1384.9999
$ws.Range("N61").Value = -4105

$ws.Range("H131").Value = 2797.0312
$ws.Range("I131").Value = 4198.75
$ws.Range("J131").Value = 2329.7917
$ws.Range("K131").Value = 12596.25
$ws.Range("L131").Value = 6989.375100000001
$ws.Range("M131").Value = -7556.25
$ws.Range("N131").Value = -17069.3751

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5059.8887
$ws.Range("I70").Value = 4456.5835
$ws.Range("J70").Value = 6266.5
$ws.Range("K70").Value = 4456.5835
$ws.Range("L70").Value = 6266.5
$ws.Range("M70").Value = -4186.5835
$ws.Range("N70").Value = -6806.5

$ws.Range("H73").Value = 5059.8887
$ws.Range("I73").Value = 4456.5835
$ws.Range("J73").Value = 6266.5
$ws.Range("K73").Value = 4456.5835
$ws.Range("L73").Value = 6266.5
$ws.Range("M73").Value = -3520.5835
$ws.Range("N73").Value = -8138.5

$ws.Range("H80").Value = 4509.231
$ws.Range("I80").Value = 5764
$ws.Range("J80").Value = 2798.182
$ws.Range("K80").Value = 5764
$ws.Range("L80").Value = 2798.182
$ws.Range("M80").Value = -4766
$ws.Range("N80").Value = -4794.182

$ws.Range("H83").Value = 4509.231
$ws.Range("I83").Value = 5764
$ws.Range("J83").Value = 2798.182
$ws.Range("K83").Value = 28820
$ws.Range("L83").Value = 13990.91
$ws.Range("M83").Value = -23828
$ws.Range("N83").Value = -23974.91

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 12832543
$ws.Range("I132").Value = 5500
$ws.Range("K132").Value = 16500
$ws.Range("M132").Value = -13970

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 42162.4
$ws.Range("I122").Value = 60870.883
$ws.Range("J122").Value = 2406.875
$ws.Range("K122").Value = 182612.649
$ws.Range("L122").Value = 7220.625
$ws.Range("M122").Value = -180162.649
$ws.Range("N122").Value = -12120.625

$ws.Range("H132").Value = 1949.3572
$ws.Range("I132").Value = 1090.973
$ws.Range("J132").Value = 3620.9473
$ws.Range("K132").Value = 3272.919
$ws.Range("L132").Value = 10862.8419
$ws.Range("M132").Value = -742.9189999999999
$ws.Range("N132").Value = -15922.8419
